# Demo Verification Script fixes
# Updates the "DateProd" (column B) timestamp values on each sheet to
# reflect the latest demo-verification run, as produced by the Katalon
# test automation tooling that stamps these workbooks.

$wb = $excel.ActiveWorkbook

# --- CC-Payments-Auth: rows 2-7 -----------------------------------------
$wsAuth = $wb.Worksheets.Item("CC-Payments-Auth")
$authTimestamps = @(
    "Tue Aug 12 01:39:06 IST 2025",
    "Tue Aug 12 01:39:47 IST 2025",
    "Tue Aug 12 01:40:32 IST 2025",
    "Tue Aug 12 01:41:12 IST 2025",
    "Tue Aug 12 01:41:47 IST 2025",
    "Tue Aug 12 01:42:26 IST 2025"
)
for ($i = 0; $i -lt $authTimestamps.Length; $i++) {
    $row = $i + 2
    $wsAuth.Cells.Item($row, 2).Value = $authTimestamps[$i]
}

# --- ACH-Payments-Debit: rows 2-10 --------------------------------------
$wsDebit = $wb.Worksheets.Item("ACH-Payments-Debit")
$debitTimestamps = @(
    "Tue Aug 12 01:43:04 IST 2025",
    "Tue Aug 12 01:43:44 IST 2025",
    "Tue Aug 12 01:44:23 IST 2025",
    "Tue Aug 12 01:45:03 IST 2025",
    "Tue Aug 12 01:45:42 IST 2025",
    "Tue Aug 12 01:46:21 IST 2025",
    "Tue Aug 12 01:47:01 IST 2025",
    "Tue Aug 12 01:47:43 IST 2025",
    "Tue Aug 12 01:48:23 IST 2025"
)
for ($i = 0; $i -lt $debitTimestamps.Length; $i++) {
    $row = $i + 2
    $wsDebit.Cells.Item($row, 2).Value = $debitTimestamps[$i]
}

# --- CC-Payments-Sale: row 2 --------------------------------------------
$wsSale = $wb.Worksheets.Item("CC-Payments-Sale")
$wsSale.Cells.Item(2, 2).Value = "Tue Aug 12 01:49:04 IST 2025"
